# Append: 2025-12-18 01:21 JST
# Update the "取得日時" (retrieved-at) timestamp in column A for all existing
# data rows (2-14) on the "ランサーズ" sheet to the new run timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-18 01:21:20"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
